$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws 'D2' '42.810.07'
Set-TextCell $ws 'E2' '  +0.58%  '

Set-TextCell $ws 'D3' '2.306.18'
Set-TextCell $ws 'E3' '  +0.08%  '

Set-TextCell $ws 'E4' '  -0.03%  '

Set-TextCell $ws 'D5' '318.08'
Set-TextCell $ws 'E5' '  -0.03%  '

Set-TextCell $ws 'D6' '104.52'
Set-TextCell $ws 'E6' '  +1.18%  '

Set-TextCell $ws 'D7' '0.629'
Set-TextCell $ws 'E7' '  -0.43%  '

Set-TextCell $ws 'E8' '  -0.12%  '

Set-TextCell $ws 'D9' '0.605'
Set-TextCell $ws 'E9' '  -0.94%  '

Set-TextCell $ws 'D10' '39.90'
Set-TextCell $ws 'E10' '  +0.56%  '

Set-TextCell $ws 'D11' '0.0910'
Set-TextCell $ws 'E11' '  -0.15%  '

Set-TextCell $ws 'D12' '8.52'
Set-TextCell $ws 'E12' '  +1.93%  '

Set-TextCell $ws 'E13' '  +2.54%  '

Set-TextCell $ws 'E14' '  +4.57%  '

Set-TextCell $ws 'D15' '15.44'
Set-TextCell $ws 'E15' '  +0.73%  '

Set-TextCell $ws 'D16' '2.655.89'
Set-TextCell $ws 'E16' '  +0.14%  '

Set-TextCell $ws 'D17' '2.306.28'
Set-TextCell $ws 'E17' '  -0.16%  '

Set-TextCell $ws 'D18' '42.739.69'
Set-TextCell $ws 'E18' '  +0.57%  '

Set-TextCell $ws 'D19' '14.83'
Set-TextCell $ws 'E19' '  +35.99%  '

Set-TextCell $ws 'D20' '7.56'
Set-TextCell $ws 'E20' '  +0.76%  '

Set-TextCell $ws 'E21' '  +0.27%  '

Set-TextCell $ws 'D22' '73.98'
Set-TextCell $ws 'E22' '  +1.06%  '

Set-TextCell $ws 'D23' '3.57'
Set-TextCell $ws 'E23' '  -0.63%  '

Set-TextCell $ws 'D24' '267.10'
Set-TextCell $ws 'E24' '  -5.22%  '

Set-TextCell $ws 'E25' '  -1.66%  '

Set-TextCell $ws 'E26' '  +0.46%  '

Set-TextCell $ws 'D27' '10.99'
Set-TextCell $ws 'E27' '  +0.98%  '

Set-TextCell $ws 'D28' '2.35'
Set-TextCell $ws 'E28' '  +0.06%  '

Set-TextCell $ws 'D29' '6.77'
Set-TextCell $ws 'E29' '  +14.40%  '

Set-TextCell $ws 'D30' '22.67'
Set-TextCell $ws 'E30' '  -1.13%  '

Set-TextCell $ws 'D31' '37.57'
Set-TextCell $ws 'E31' '  +4.16%  '

Set-TextCell $ws 'D32' '166.18'
Set-TextCell $ws 'E32' '  +0.98%  '

Set-TextCell $ws 'D33' '0.0885'
Set-TextCell $ws 'E33' '  +1.11%  '

Set-TextCell $ws 'E34' '  -3.91%  '

Set-TextCell $ws 'D35' '2.61'
Set-TextCell $ws 'E35' '  -0.43%  '

Set-TextCell $ws 'E36' '  -3.33%  '

Set-TextCell $ws 'D37' '4.58'
Set-TextCell $ws 'E37' '  -0.67%  '

Set-TextCell $ws 'E38' '  -3.98%  '

Set-TextCell $ws 'D39' '3.73'
Set-TextCell $ws 'E39' '  +0.06%  '

Set-TextCell $ws 'D40' '2.72'
Set-TextCell $ws 'E40' '  -2.57%  '

Set-TextCell $ws 'D41' '1.60'
Set-TextCell $ws 'E41' '  +6.73%  '

Set-TextCell $ws 'D42' '70.57'
Set-TextCell $ws 'E42' '  +1.16%  '

Set-TextCell $ws 'B43' 'Algorand'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws 'D43' '0.230'
Set-TextCell $ws 'E43' '  +1.45%  '

Set-TextCell $ws 'B44' 'BitcoinSV'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextCell $ws 'D44' '95.62'
Set-TextCell $ws 'E44' '  -3.57%  '

Set-TextCell $ws 'E45' '  -0.03%  '

Set-TextCell $ws 'D46' '12.29'
Set-TextCell $ws 'E46' '  +1.14%  '

Set-TextCell $ws 'B47' 'Aave'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D47' '115.96'
Set-TextCell $ws 'E47' '  +3.06%  '

Set-TextCell $ws 'B48' 'ordi'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextCell $ws 'D48' '81.46'
Set-TextCell $ws 'E48' '  +3.18%  '

Set-TextCell $ws 'D49' '1.700.72'
Set-TextCell $ws 'E49' '  +5.62%  '

Set-TextCell $ws 'E50' '  -0.90%  '

Set-TextCell $ws 'E51' '  -2.20%  '
